# "Added planning to final version"
# Fill in the "Uitgevoerde uren" (hours actually worked) columns (D, F, H, J)
# for each team member (Jip/Tim/Leo/Rik) across both weekly tables, and add
# the missing per-person week-2 totals formulas in row 42.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Week 1 table (rows 4-17) ---
$ws.Range("D4").Value = 6
$ws.Range("F4").Value = 8
$ws.Range("H4").Value = 4
$ws.Range("J4").Value = 2

$ws.Range("D6").Value = 5
$ws.Range("H6").Value = 5
$ws.Range("J6").Value = 2

$ws.Range("F7").Value = 2

$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0

$ws.Range("D10").Value = 2
$ws.Range("J10").Value = 7

$ws.Range("D11").Value = 2
$ws.Range("F11").Value = 2

$ws.Range("F12").Value = 2
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0

$ws.Range("F13").Value = 3
$ws.Range("J13").Value = 4

$ws.Range("J15").Value = 1

$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = 6
$ws.Range("J16").Value = 3

$ws.Range("D17").Value = 3

# --- Week 2 table (rows 22-39) ---
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 6
$ws.Range("J22").Value = 2

$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0

$ws.Range("D24").Value = 7
$ws.Range("F24").Value = 4
$ws.Range("H24").Value = 6
$ws.Range("J24").Value = 16

$ws.Range("D25").Value = 2
$ws.Range("F25").Value = 4
$ws.Range("H25").Value = 2
$ws.Range("J25").Value = 2

$ws.Range("D26").Value = 4
$ws.Range("H26").Value = 20
$ws.Range("J26").Value = 6

$ws.Range("J27").Value = 0.5

$ws.Range("D31").Value = 3
$ws.Range("J31").Value = 2

$ws.Range("D32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0

$ws.Range("D33").Value = 4
$ws.Range("F33").Value = 1
$ws.Range("H33").Value = 1
$ws.Range("J33").Value = 1

$ws.Range("D34").Value = 1
$ws.Range("F34").Value = 1

$ws.Range("H35").Value = 2.5
$ws.Range("J35").Value = 2.5

$ws.Range("F36").Value = 2
$ws.Range("H36").Value = 2

$ws.Range("D37").Value = 2
$ws.Range("F37").Value = 0
$ws.Range("H37").Value = 0

$ws.Range("D38").Value = 16
$ws.Range("F38").Value = 16
$ws.Range("H38").Value = 6
$ws.Range("J38").Value = 10

# --- Grand-total row (42): add the per-person "Uitgevoerde uren" sums that
# were still missing ---
$ws.Range("D42").Formula = "=SUM(D40,D29,D19)"
$ws.Range("F42").Formula = "=SUM(F29,F40,F19)"
$ws.Range("H42").Formula = "=SUM(H40,H29,H19)"
$ws.Range("J42").Formula = "=SUM(J40,J29,J19)"

# --- View state: scroll position / active selection as left by the author ---
[void]$ws.Range("G39").Select()
$excel.ActiveWindow.ScrollRow = 13
